# Income sheet refactor: drop the "Interest From Saving Account" entry and
# refresh the remaining income sources with new amounts/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Income")

# Remove row 2 ("Interest From Saving Account") entirely; Excel shifts the
# rows below it up by one (rows 3-6 become rows 2-5).
$ws.Rows.Item(2).Delete()

# Row 2 -> Trading
$ws.Range("A2").Value = "Trading"
$ws.Range("B2").Value = 25000
$ws.Range("C2").Value = 45853.22928240741

# Row 3 -> Business Income
$ws.Range("A3").Value = "Business Income"
$ws.Range("B3").Value = 65000
$ws.Range("C3").Value = 45849.22928240741

# Row 4 -> Youtube Revenue
$ws.Range("A4").Value = "Youtube Revenue"
$ws.Range("B4").Value = 35000
$ws.Range("C4").Value = 45848.22928240741

# Row 5 -> Salary
$ws.Range("A5").Value = "Salary"
$ws.Range("B5").Value = 40000
$ws.Range("C5").Value = 45839.22928240741
